$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Q8" in J1, matching the style of the existing headers (copy format from I1)
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Q8"

# Update simulated rt_data values across the grid (bugfixed evaluation values)
$ws.Range("B2").Value = 3.534544165640355
$ws.Range("C2").Value = 0.6043611617905069
$ws.Range("D2").Value = 0.3677908113261881
$ws.Range("E2").Value = 0.6933852433072614
$ws.Range("F2").Value = 1.257110715986727
$ws.Range("G2").Value = 0.9709551209826655
$ws.Range("H2").Value = 1.092800686503665
$ws.Range("B3").Value = 0.1079662819227423
$ws.Range("C3").Value = -0.1286040685415766
$ws.Range("D3").Value = 0.1969903634394967
$ws.Range("E3").Value = 0.7607158361189619
$ws.Range("F3").Value = 0.4745602411149009
$ws.Range("G3").Value = 0.5964058066359001
$ws.Range("B4").Value = -0.6261225043750751
$ws.Range("C4").Value = -0.3005280723940018
$ws.Range("D4").Value = 0.2631974002854633
$ws.Range("E4").Value = -0.02295819471859767
$ws.Range("F4").Value = 0.09888737080240162
$ws.Range("G4").Value = -0.3347921734358863
$ws.Range("H4").Value = 1.551084065451221
$ws.Range("I4").Value = -0.1524750000452286
$ws.Range("J4").Value = -0.2080565315694543
$ws.Range("B5").Value = 0.294793270082792
$ws.Range("C5").Value = 0.8585187427622571
$ws.Range("D5").Value = 0.5723631477581961
$ws.Range("E5").Value = 0.6942087132791954
$ws.Range("F5").Value = 0.2605291690409075
$ws.Range("G5").Value = 2.146405407928015
$ws.Range("H5").Value = 0.4428463424315652
$ws.Range("I5").Value = 0.3872648109073396
$ws.Range("B6").Value = 0.4943829294508504
$ws.Range("C6").Value = 0.2082273344467894
$ws.Range("D6").Value = 0.3300728999677887
$ws.Range("E6").Value = -0.1036066442704993
$ws.Range("F6").Value = 1.782269594616608
$ws.Range("G6").Value = 0.07871052912015841
$ws.Range("H6").Value = 0.02312899759593279
$ws.Range("B7").Value = -0.2133677787764014
$ws.Range("C7").Value = -0.09152221325540211
$ws.Range("D7").Value = -0.52520175749369
$ws.Range("E7").Value = 1.360674481393417
$ws.Range("F7").Value = -0.3428845841030324
$ws.Range("G7").Value = -0.398466115627258
$ws.Range("B8").Value = 0.0579608153039004
$ws.Range("C8").Value = -0.3757187289343875
$ws.Range("D8").Value = 1.51015750995272
$ws.Range("E8").Value = -0.1934015555437298
$ws.Range("F8").Value = -0.2489830870679555
$ws.Range("G8").Value = 0.5379667755061348
$ws.Range("H8").Value = -1.651730635291429
$ws.Range("I8").Value = -0.8137395822194835
$ws.Range("B9").Value = -0.4295258376674695
$ws.Range("C9").Value = 1.456350401219638
$ws.Range("D9").Value = -0.2472086642768118
$ws.Range("E9").Value = -0.3027901958010374
$ws.Range("F9").Value = 0.4841596667730528
$ws.Range("G9").Value = -1.705537744024511
$ws.Range("H9").Value = -0.8675466909525655
$ws.Range("B10").Value = 1.781099436349905
$ws.Range("C10").Value = 0.07754037085345544
$ws.Range("D10").Value = 0.02195883932922982
$ws.Range("E10").Value = 0.8089087019033201
$ws.Range("F10").Value = -1.380788708894244
$ws.Range("G10").Value = -0.5427976558222982
$ws.Range("B11").Value = -0.2578266341250811
$ws.Range("C11").Value = -0.3134081656493067
$ws.Range("D11").Value = 0.4735416969247836
$ws.Range("E11").Value = -1.716155713872781
$ws.Range("F11").Value = -0.8781646608008347
$ws.Range("B12").Value = -0.4069492870295619
$ws.Range("C12").Value = 0.3800005755445284
$ws.Range("D12").Value = -1.809696835253036
$ws.Range("E12").Value = -0.9717057821810899
$ws.Range("B13").Value = 0.8628474748582136
$ws.Range("C13").Value = -1.326849935939351
$ws.Range("D13").Value = -0.4888588828674048
$ws.Range("B14").Value = -1.436438518536832
$ws.Range("C14").Value = -0.5984474654648866
$ws.Range("B15").Value = -0.6276069079710285
